$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "cards" / "Visa" column values and document-weight scores to row 2
$ws.Range("D2").Value = "cards"
$ws.Range("E2").Value = "Visa"
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 7
$ws.Range("M2").Value = 7

# Match the author's on-screen selection/scroll position when the file was saved
$ws.Range("L2").Select()
